$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Payment_Method sheet: add the PayPal locator rows (A4:C13)
# ---------------------------------------------------------------------------
$wsPay = $wb.Worksheets.Item("Payment_Method")

# A4/A5 already carry the "s=5" locator-name style; stamp the same style
# onto the brand-new A6:A13 cells before filling them in.
$wsPay.Range("A5").Copy($wsPay.Range("A6:A13"))

$wsPay.Range("A4").Value = "paypal email input"
$wsPay.Range("B4").Value = "id"
$wsPay.Range("C4").Value = "email"

$wsPay.Range("A5").Value = "paypal email input"
$wsPay.Range("B5").Value = "name"
$wsPay.Range("C5").Value = "login_email"

$wsPay.Range("A6").Value = "paypal next button"
$wsPay.Range("B6").Value = "id"
$wsPay.Range("C6").Value = "btnNext"

$wsPay.Range("A7").Value = "paypal next button"
$wsPay.Range("B7").Value = "name"
$wsPay.Range("C7").Value = "btnNext"

$wsPay.Range("A8").Value = "paypal password"
$wsPay.Range("B8").Value = "id"
$wsPay.Range("C8").Value = "password"

$wsPay.Range("A9").Value = "paypal password"
$wsPay.Range("B9").Value = "name"
$wsPay.Range("C9").Value = "login-password"

$wsPay.Range("A10").Value = "paypal login button"
$wsPay.Range("B10").Value = "id"
$wsPay.Range("C10").Value = "btnLogin"

$wsPay.Range("A11").Value = "paypal login button"
$wsPay.Range("B11").Value = "name"
$wsPay.Range("C11").Value = "btnLogin"

$wsPay.Range("A12").Value = "paypal pay now button"
$wsPay.Range("B12").Value = "id"
$wsPay.Range("C12").Value = "payment-submit-btn"

$wsPay.Range("A13").Value = "paypal pay now button"
$wsPay.Range("B13").Value = "xpath"
$wsPay.Range("C13").Value = '//*[@id="payment-submit-btn"]'

$wsPay.Activate()
$wsPay.Range("A12").Select()

# ---------------------------------------------------------------------------
# Checkout_Page sheet: add the "order number" locator row, tidy up the
# trailing blank rows.
# ---------------------------------------------------------------------------
$wsChk = $wb.Worksheets.Item("Checkout_Page")

# row 15 becomes completely empty (drop the lone formatted D15 cell)
$wsChk.Range("D15").Clear()

# row 16 gains the "order number" xpath locator; B16/C16 need to pick up
# the same style already used by the other xpath rows (e.g. B14/C14).
$wsChk.Range("A16").Value = "order number"

$wsChk.Range("C14").Copy($wsChk.Range("B16"))
$wsChk.Range("B16").Value = "xpath"

$wsChk.Range("C14").Copy($wsChk.Range("C16"))
$wsChk.Range("C16").Value = '//*[@id="__next"]/div/div[2]/div[2]/div[2]/div[1]/div/div/span'

# row 19 gets a blank, formatted D cell (matching D17/D18's style)
$wsChk.Range("D18").Copy($wsChk.Range("D19"))

# a new trailing blank row 22 (matching A20/A21's style)
$wsChk.Range("A21").Copy($wsChk.Range("A22"))

$wsChk.Activate()
$wsChk.Range("A16").Select()

# Leave Payment_Method as the active sheet/tab (matches the workbook's
# unchanged activeTab=9 / tabSelected on sheet10).
$wsPay.Activate()
